$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Starting layout (before):
#   row 12 = last existing item row (item 6)
#   row 13 = totals row (P13:Q13 merged, value 486.42)
#   row 14 = footer row (A14:F14 timestamp, G14:I14 page, K14:Q14 credit)
#
# Target layout (after):
#   row 12 = unchanged (item 6)
#   row 13 = NEW item row (item 7: TOBRIN 0.3% EYE DROPS 5 ML)
#   row 14 = totals row (moved down from 13, new total 509.42)
#   row 15 = footer row (moved down from 14, timestamp updated to 10:12 AM)
# ---------------------------------------------------------------------------

# 1) Unmerge everything that is about to move / be rebuilt
$ws.Range("P13:Q13").UnMerge()
$ws.Range("A14:F14").UnMerge()
$ws.Range("G14:I14").UnMerge()
$ws.Range("K14:Q14").UnMerge()

# 2) Move the footer row (old row 14) down to row 15
$ws.Range("A14:Q14").Copy()
$ws.Range("A15:Q15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Wednesday, 18 June, 2025 10:12 AM"
$ws.Range("G15").Value = "1/1"
$ws.Range("K15").Value = "developed by : Abdelaziz Talaat"
$ws.Rows.Item(15).RowHeight = 16.5

# clear the old row 14 cells that are not going to be overwritten below (D14,E14,F14 etc.)
$ws.Range("A14:Q14").ClearContents()
$ws.Range("A14:Q14").ClearFormats()

# 3) Move the totals row (old row 13) down to row 14, with the updated sum
$ws.Range("P13:Q13").Copy()
$ws.Range("P14:Q14").PasteSpecial(-4122)
$ws.Range("P14").Value = 509.42
$ws.Rows.Item(14).RowHeight = 25.5

$ws.Range("P13:Q13").ClearContents()
$ws.Range("P13:Q13").ClearFormats()

# 4) Build the new item row 13, copying the look of row 12 (same style pattern)
$ws.Range("A12:Q12").Copy()
$ws.Range("A13:Q13").PasteSpecial(-4122)

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "TOBRIN 0.3% EYE DROPS 5 ML"
$ws.Range("H13").Value = "4:0"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "23.00"
$ws.Range("P13").Value = "23.0000"
$ws.Range("Q13").Value = "1:0"
$ws.Rows.Item(13).RowHeight = 24.75

# 5) Re-create merged cells at their new locations
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()
$ws.Range("P14:Q14").Merge()
$ws.Range("A15:F15").Merge()
$ws.Range("G15:I15").Merge()
$ws.Range("K15:Q15").Merge()

Write-Output "edit complete"
